$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 12 (VatanseverSakin2009267 entry) entirely; rows below shift up
$ws.Rows.Item(12).Delete()

# Mark row 10 (Manenti201840) as extracted
$ws.Range("N10").Value = "y"

# Complete extraction for row 11 (Peng201432)
$ws.Range("M11").Value = "acclimation; excluded the figures/data that reported the fluctuating treatments at the two different min/max. "
$ws.Range("K11").Value = "fig 1"
$ws.Range("N11").Value = "y"

# Rebuild the autofilter range to reflect the removed row
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A2:T19").AutoFilter()
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$2:`$T`$19"

# Rebuild the remembered sort state to reflect the removed row
$sortRange = $ws.Range("A2:T17")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("J2:J17"), 0, 2, 0, 0)
$ws.Sort.SortFields.Add($ws.Range("O2:O17"), 0, 2, 0, 0)
$ws.Sort.SortFields.Add($ws.Range("A2:A17"), 0, 1, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Restore selection as saved by the author
$ws.Range("D20").Select()

# Restore the window position as saved by the author
$win = $excel.ActiveWindow
$win.Left = 2840
$win.Top = 460
